# Updated cryptos list (price + 1h volume change %) pulled on
# Fri Sep  1 18:56:12 UTC 2023 with GitHub Actions.
# Cells D2:E51 hold price/volume text (and a handful of coins in rows
# 36-38 / 41-42 re-ranked, shifting B/C/D/E down a row). Values are
# written with a leading "'" via Value2 so Excel keeps them as literal
# text (several, e.g. "25.992.86" or "1.009", would otherwise be
# auto-parsed as a date/number) and ClearFormats() afterwards strips the
# resulting quote-prefix cell style back off so formatting is untouched.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value2 = "'25.992.86"
$ws.Range("E2").Value2 = "'  -1.73%  "
$ws.Range("D3").Value2 = "'1.637.11"
$ws.Range("E3").Value2 = "'  -1.88%  "
$ws.Range("D4").Value2 = "'1.009"
$ws.Range("E4").Value2 = "'  +0.18%  "
$ws.Range("D5").Value2 = "'215.45"
$ws.Range("E5").Value2 = "'  -1.49%  "
$ws.Range("D6").Value2 = "'0.5022"
$ws.Range("E6").Value2 = "'  -2.55%  "
$ws.Range("D7").Value2 = "'1.011"
$ws.Range("E7").Value2 = "'  +0.44%  "
$ws.Range("D8").Value2 = "'0.2577"
$ws.Range("D9").Value2 = "'0.06417"
$ws.Range("E9").Value2 = "'  -0.81%  "
$ws.Range("D10").Value2 = "'19.53"
$ws.Range("E10").Value2 = "'  -2.44%  "
$ws.Range("D11").Value2 = "'0.07744"
$ws.Range("E11").Value2 = "'  +0.99%  "
$ws.Range("D12").Value2 = "'1.645.72"
$ws.Range("E12").Value2 = "'  -1.47%  "
$ws.Range("D13").Value2 = "'4.260"
$ws.Range("E13").Value2 = "'  -2.04%  "
$ws.Range("D14").Value2 = "'1.861.15"
$ws.Range("E14").Value2 = "'  -1.94%  "
$ws.Range("D15").Value2 = "'0.5453"
$ws.Range("E15").Value2 = "'  -2.00%  "
$ws.Range("D16").Value2 = "'0.0₅7969"
$ws.Range("E16").Value2 = "'  -1.14%  "
$ws.Range("D17").Value2 = "'63.52"
$ws.Range("E17").Value2 = "'  -1.96%  "
$ws.Range("D18").Value2 = "'25.992.51"
$ws.Range("E18").Value2 = "'  -1.82%  "
$ws.Range("D19").Value2 = "'1.011"
$ws.Range("E19").Value2 = "'  +0.41%  "
$ws.Range("D20").Value2 = "'205.35"
$ws.Range("E20").Value2 = "'  -2.40%  "
$ws.Range("D21").Value2 = "'4.310"
$ws.Range("E21").Value2 = "'  -2.69%  "
$ws.Range("D22").Value2 = "'10.00"
$ws.Range("E22").Value2 = "'  -1.24%  "
$ws.Range("D23").Value2 = "'5.974"
$ws.Range("E23").Value2 = "'  +1.31%  "
$ws.Range("D24").Value2 = "'1.011"
$ws.Range("E24").Value2 = "'  +0.38%  "
$ws.Range("D25").Value2 = "'1.960"
$ws.Range("E25").Value2 = "'  +12.78%  "
$ws.Range("D26").Value2 = "'141.57"
$ws.Range("E26").Value2 = "'  -2.49%  "
$ws.Range("D27").Value2 = "'0.1152"
$ws.Range("E27").Value2 = "'  -1.02%  "
$ws.Range("D28").Value2 = "'15.79"
$ws.Range("E28").Value2 = "'  -0.03%  "
$ws.Range("D29").Value2 = "'6.813"
$ws.Range("E29").Value2 = "'  -2.86%  "
$ws.Range("D30").Value2 = "'1.238"
$ws.Range("E30").Value2 = "'  -1.87%  "
$ws.Range("D31").Value2 = "'0.05004"
$ws.Range("E31").Value2 = "'  -4.14%  "
$ws.Range("D32").Value2 = "'3.272"
$ws.Range("E32").Value2 = "'  -2.99%  "
$ws.Range("E33").Value2 = "'  -0.78%  "
$ws.Range("D34").Value2 = "'1.540"
$ws.Range("E34").Value2 = "'  -2.94%  "
$ws.Range("D35").Value2 = "'2.339"
$ws.Range("E35").Value2 = "'  -1.62%  "
$ws.Range("B36").Value2 = "'ImmutableX"
$ws.Range("C36").Value2 = "'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D36").Value2 = "'0.5669"
$ws.Range("E36").Value2 = "'  -1.34%  "
$ws.Range("B37").Value2 = "'ARBITRUM"
$ws.Range("C37").Value2 = "'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D37").Value2 = "'0.8892"
$ws.Range("E37").Value2 = "'  -3.82%  "
$ws.Range("B38").Value2 = "'MXToken"
$ws.Range("C38").Value2 = "'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D38").Value2 = "'2.606"
$ws.Range("E38").Value2 = "'  -5.50%  "
$ws.Range("D39").Value2 = "'1.119.81"
$ws.Range("E39").Value2 = "'  -3.96%  "
$ws.Range("D40").Value2 = "'0.01566"
$ws.Range("E40").Value2 = "'  -2.12%  "
$ws.Range("B41").Value2 = "'PaxDollar"
$ws.Range("C41").Value2 = "'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D41").Value2 = "'1.011"
$ws.Range("E41").Value2 = "'  +0.42%  "
$ws.Range("B42").Value2 = "'mCoin"
$ws.Range("C42").Value2 = "'https://coinranking.com/coin/fzVgyjBcRc9+mcoin-mcoin"
$ws.Range("D42").Value2 = "'2.568"
$ws.Range("E42").Value2 = "'  -0.21%  "
$ws.Range("D43").Value2 = "'5.617"
$ws.Range("E43").Value2 = "'  -0.53%  "
$ws.Range("D44").Value2 = "'0.8165"
$ws.Range("E44").Value2 = "'  -3.01%  "
$ws.Range("D45").Value2 = "'99.75"
$ws.Range("E45").Value2 = "'  -0.50%  "
$ws.Range("D46").Value2 = "'1.771.28"
$ws.Range("E46").Value2 = "'  -2.05%  "
$ws.Range("D47").Value2 = "'0.0₈111"
$ws.Range("E47").Value2 = "'  -0.62%  "
$ws.Range("D48").Value2 = "'0.4536"
$ws.Range("E48").Value2 = "'  +0.79%  "
$ws.Range("D50").Value2 = "'54.73"
$ws.Range("E50").Value2 = "'  -2.50%  "
$ws.Range("D51").Value2 = "'0.05037"
$ws.Range("E51").Value2 = "'  -1.61%  "

$ws.Range("B2:E51").ClearFormats()
